$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Update row 4 (D4, F4 changed; G4 added) ---
$ws.Range("D4").Value = 1
$ws.Range("F4").Value = 44203

$ws.Range("E4").Copy()
$ws.Range("G4").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("G4").Value = 44203

# --- Update row 5 (D5, E5, F5 changed) ---
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 44204
$ws.Range("F5").Value = 44209

# --- Add row 6 (new) ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Documentation for JSON data"
$ws.Range("C6").Value = "LMS v3"
$ws.Range("D6").Value = 1

$ws.Range("E4:F4").Copy()
$ws.Range("E6:F6").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("E6").Value = 44210
$ws.Range("F6").Value = 44210

# --- Add row 7 (new) ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Documentation for existing/planned features"
$ws.Range("C7").Value = "LMS v3"

# --- Update selection ---
$ws.Range("H10").Select()
